# Update the vm_pu results for the "380 kV" case (B column slack voltage
# changed from 1.05 to 1.02 pu, which cascades into recomputed bus voltages).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns (in sheet order) that receive new values: B,C,D,E,F,I,J,K,L,M,N
$cols = @("B","C","D","E","F","I","J","K","L","M","N")

# New values per data row (rows 2..25 in the worksheet), one array per row
# ordered to match $cols above.
$data = @{
    2  = @("1.02","1.050339374274593","1.060295203387247","1.057623515381107","1.068828435398397","1.039745538513281","1.055373050793867","1.063022372017043","1.060357977551843","1.071532573017533","1.056871801578976")
    3  = @("1.02","1.051356741733243","1.06124227267802","1.058526189964492","1.069808718031774","1.039896116167578","1.05603972913444","1.06378350363224","1.061074299364995","1.072328500920605","1.057539426679261")
    4  = @("1.02","1.052015485484883","1.061855801817626","1.059111013771675","1.070443812989547","1.039992213332817","1.056470939909139","1.064276091559682","1.061537904540922","1.072843683869005","1.057971249822721")
    5  = @("1.02","1.052292525816865","1.062113898861383","1.059357048097643","1.070710994113517","1.040032292148273","1.056652178313113","1.064483194958416","1.061732826530451","1.073060305184848","1.058152745606017")
    6  = @("1.02","1.052339048214868","1.062157244393372","1.059398368538895","1.070755865945769","1.040039002754727","1.056682606517447","1.064517969644276","1.061765556077113","1.073096679081519","1.058183217021891")
    7  = @("1.02","1.052019186901688","1.061859249858678","1.059114300610473","1.070447382342592","1.039992750127536","1.056473361793957","1.064278858810179","1.061540509011256","1.072846578223913","1.057973675146892")
    8  = @("1.02","1.050683107095678","1.060615121832865","1.057928426061209","1.069159562985077","1.039796703434901","1.055598393559507","1.063279581902795","1.060600041245158","1.071801525974743","1.057097464357189")
    9  = @("1.02","1.048332153647219","1.058428308733631","1.055844425121214","1.066896333053907","1.039441026278743","1.054055284347665","1.06151941738458","1.058943598316139","1.069961311004462","1.05555216375327")
    10 = @("1.02","1.046767166132669","1.056974192252217","1.054458960764983","1.065391659600817","1.0391970638362","1.053025713974822","1.060346489056974","1.057839878609618","1.068735430520118","1.054521131272404")
    11 = @("1.02","1.046090065179209","1.056345447849456","1.053859970148562","1.064741115904931","1.039089807250618","1.052579710764256","1.059838730662506","1.057362102533075","1.068204842411483","1.054074494686122")
    12 = @("1.02","1.04583864268944","1.056112040002756","1.05363761827277","1.064499624724162","1.039049724493041","1.052414017092021","1.059650146505072","1.057184657312699","1.068007793205554","1.053908565709867")
    13 = @("1.02","1.045892569905506","1.056162100616642","1.053685307157672","1.064551418611259","1.0390583333655","1.052449560213144","1.059690597532944","1.057222718902819","1.068050059338728","1.053944159306299")
    14 = @("1.02","1.046069280822798","1.0563261515078","1.053841587613255","1.06472115110869","1.039086498949937","1.052566015045328","1.059823141826316","1.057347434400725","1.068188553546942","1.054060779517701")
    15 = @("1.02","1.046178169257457","1.056427246702353","1.053937895729096","1.064825748770693","1.039103820522902","1.052637762992011","1.05990480940303","1.057424278685203","1.06827388904198","1.054132629354699")
    16 = @("1.02","1.046812114699604","1.057015938941488","1.054498733364871","1.065434854983586","1.039204148019986","1.053055309743389","1.060380190059144","1.057871590089043","1.068770648744666","1.054550769070357")
    17 = @("1.02","1.04720991923562","1.057385451028651","1.054850780241754","1.065817196932823","1.039266647564467","1.053317174753516","1.060678418428228","1.05815221523797","1.069082314152335","1.054813005958835")
    18 = @("1.02","1.047442004826374","1.057601067721155","1.055056212371561","1.066040305960934","1.039302946232396","1.053469897568918","1.060852382220859","1.058315912772798","1.069264125082318","1.054965945658152")
    19 = @("1.02","1.047521148979514","1.057674602091375","1.055126274508729","1.066116396531341","1.03931529662416","1.0535219689239","1.060911701428993","1.058371731690653","1.069326121577767","1.055018090960431")
    20 = @("1.02","1.047167233094845","1.057345796920119","1.054812999715895","1.065776165381256","1.039259958114628","1.053289081036556","1.060646420075808","1.05812210539949","1.069048873145616","1.054784872345575")
    21 = @("1.02","1.046017241568321","1.056277838852213","1.053795563044062","1.064671164977797","1.039078211589067","1.052531722773303","1.059784110279295","1.057310708179071","1.068147769485826","1.054026438546717")
    22 = @("1.02","1.045294676428478","1.055607157513416","1.053156670193719","1.063977274059094","1.03896253526306","1.052055378636151","1.059242057518892","1.056800678719453","1.067581411859419","1.053549417946253")
    23 = @("1.02","1.045677676204999","1.055962623528808","1.053495282272345","1.064345036276716","1.039023990527443","1.052307912889329","1.0595293987178","1.057071042576335","1.067881629319445","1.053802310827036")
    24 = @("1.02","1.047186520960855","1.057363714637519","1.054830070833916","1.065794705478395","1.039262981271902","1.053301775435509","1.060660878716156","1.058135710697349","1.069063983629873","1.05479758477203")
    25 = @("1.02","1.048939526094505","1.058992994303299","1.056382511896823","1.067480706069324","1.039534185586001","1.054454364849874","1.061974375154952","1.059371730746721","1.070436890888573","1.055951810995566")
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range("$($cols[$i])$row").Value = [double]$vals[$i]
    }
}
